$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old SUM formula in C29 and extend the table down to row 53.
# Column A / B get 0..24, column C gets 0 for every new row (and C29 loses
# its old SUM formula, becoming a literal 0 like the rest).
for ($i = 0; $i -le 24; $i++) {
    $row = 29 + $i
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $i
    $ws.Cells.Item($row, 3).Value = 0
}

# Columns A and B on the new rows pick up a dedicated font color style.
$ws.Range("A29:B53").Font.Color = 0

# Match the saved selection state.
$ws.Range("C29:C53").Select()

Write-Output "done"
